# Edit script implementing the diff changes for COVID.docx
$d = $word.ActiveDocument
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# ---------------------------------------------------------------------------
# Edit A: paragraph "Este proyecto va a dirigido..." (paragraph 3)
#   1) Split "...etc.)." into several runs incl. proofErr spellStart/spellEnd,
#      append new sentence about remote diagnostics, then break into a new
#      paragraph "El proyecto consta de 3 secciones fundamentales:"
# ---------------------------------------------------------------------------

# 1a. Replace the run containing "etc.). Para que las personas puedan llevar
#     donativos. Este proyecto consta de " -- first locate the whole
#     paragraph and find the split point using Find on the full paragraph
#     text, then use InsertParagraphAfter to split, and InsertXML to set
#     exact run/proofErr structure for each side.

$rng = $d.Content
$rng.Find.Execute("de los lugares donde exista la producción y/o distribución de productos de higiene básica, así como de materiales necesarios para hacer frente a la contingencia (mascarillas, gel antibacterial, cubre bocas, etc.). Para que las personas puedan llevar donativos. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "A1 found:" $rng.Text.Length
$rng.InsertParagraphAfter()

# Re-acquire paragraph objects after the split.
# Paragraph 3 now ends right after "...donativos. " and paragraph 4 begins
# with "Este proyecto consta de 3 secciones fundamentales."
$p3 = $d.Paragraphs.Item(3).Range
$frag3 = "<w:p $wns><w:pPr><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:t xml:space=`"preserve`">Este proyecto va a dirigido al ámbito de </w:t></w:r>" + `
  "<w:r><w:t>las Tecnologías de la Información</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> en la cual se tendrá un registro </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">de los lugares donde exista la producción y/o distribución de productos de higiene básica, así como de materiales necesarios para hacer frente a la contingencia (mascarillas, gel </w:t></w:r>" + `
  "<w:proofErr w:type=`"spellStart`"/>" + `
  "<w:r><w:t>antibacterial</w:t></w:r>" + `
  "<w:proofErr w:type=`"spellEnd`"/>" + `
  "<w:r><w:t>, cubre bocas, etc.)</w:t></w:r>" + `
  "<w:r><w:t>.</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> Este sistema contara con la posibilidad de realizar diagnósticos de forma remota.</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
  "</w:p>"
$res3 = $p3.InsertXML($frag3)
Write-Host "A1 result:" $res3

$p4 = $d.Paragraphs.Item(4).Range
Write-Host "A2 para4 text before:" $p4.Text
$frag4 = "<w:p $wns><w:pPr><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:t>E</w:t></w:r>" + `
  "<w:r><w:t>l</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> proyecto consta de </w:t></w:r>" + `
  "<w:r><w:t>3</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
  "<w:r><w:t>secciones</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> fundamentales</w:t></w:r>" + `
  "<w:r><w:t>:</w:t></w:r>" + `
  "</w:p>"
$res4 = $p4.InsertXML($frag4)
Write-Host "A2 result:" $res4

Write-Host "DONE_A"

# ---------------------------------------------------------------------------
# Edit C: insert a brand-new paragraph after the "Aplicación Móvil" paragraph
# (the one ending "...diagnostico a través de un chat privado."), before the
# "EJe al que corresponde" heading.
# ---------------------------------------------------------------------------
$rngC = $d.Content
$rngC.Find.Execute("diagnostico a través de un chat privado.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "C found:" $rngC.Text

# Determine the paragraph index of the found range BEFORE inserting, then the
# new empty paragraph will be at index+1 afterwards.
$idxC = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -le $rngC.Start -and $d.Paragraphs.Item($i).Range.End -ge $rngC.End) {
        $idxC = $i
    }
}
Write-Host "C paragraph index:" $idxC

$rngC.InsertParagraphAfter()

# The newly created (currently empty) paragraph is right after it.
$newParaC = $d.Paragraphs.Item($idxC + 1).Range
Write-Host "C newPara text before:" $newParaC.Text
$fragC = "<w:p $wns><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:t>La aplicación también será ca</w:t></w:r>" + `
  "<w:r><w:t>paz de lanzar notificaciones push con información</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> relevante </w:t></w:r>" + `
  "<w:r><w:t>para la población.</w:t></w:r>" + `
  "</w:p>"
$resC = $newParaC.InsertXML($fragC)
Write-Host "C result:" $resC

Write-Host "DONE_C"

# ---------------------------------------------------------------------------
# Edit D: paragraph "Que la población tenga..." -- replace entirely with new
# wording and add jc=both paragraph formatting.
# ---------------------------------------------------------------------------
$rngD = $d.Content
$rngD.Find.Execute("Que la población tenga una atención medica sin tener que salir de casa. Un sistema centralizado donde", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "D found:" $rngD.Text
$idxD = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -le $rngD.Start -and $d.Paragraphs.Item($i).Range.End -ge $rngD.End) {
        $idxD = $i
    }
}
Write-Host "D paragraph index:" $idxD
$pD = $d.Paragraphs.Item($idxD).Range
$fragD = "<w:p $wns><w:pPr><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:t xml:space=`"preserve`">Que la población tenga acceso a una atención medica tanto gratuita </w:t></w:r>" + `
  "<w:r><w:t>como</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> privada con la finalidad de evitar salir de casa</w:t></w:r>" + `
  "<w:r><w:t>, esto reducirá el riesgo de contagio por salir al exterior o visitar una institución de salud, esta característica permitirá también la optimización de recursos como el tiempo y transporte.</w:t></w:r>" + `
  "</w:p>"
$resD = $pD.InsertXML($fragD)
Write-Host "D result:" $resD

# ---------------------------------------------------------------------------
# Edit E: the final (empty) paragraph right before the sectPr -- replace
# entirely with new wording including gramStart/gramEnd proofErr markers.
# ---------------------------------------------------------------------------
$pE = $d.Paragraphs.Item($idxD + 1).Range
Write-Host "E para text before:" $pE.Text
$fragE = "<w:p $wns><w:pPr><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:t>A</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">sí mismo </w:t></w:r>" + `
  "<w:r><w:t>se podrá consultar recibir</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">un </w:t></w:r>" + `
  "<w:r><w:t>med</w:t></w:r>" + `
  "<w:r><w:t>i</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">camento </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">adecuado </w:t></w:r>" + `
  "<w:r><w:t>del lugar más cercano</w:t></w:r>" + `
  "<w:r><w:t>.</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">Ofrecerá información oficial verídica del COVID-19, esto para evitar </w:t></w:r>" + `
  "<w:proofErr w:type=`"gramStart`"/>" + `
  "<w:r><w:t>l</w:t></w:r>" + `
  "<w:r><w:t>a</w:t></w:r>" + `
  "<w:proofErr w:type=`"gramEnd`"/>" + `
  "<w:r><w:t xml:space=`"preserve`"> des información de la población y mantenerlos al día</w:t></w:r>" + `
  "<w:r><w:t>.</w:t></w:r>" + `
  "</w:p>"
$resE = $pE.InsertXML($fragE)
Write-Host "E result:" $resE

Write-Host "DONE_E"

# ---------------------------------------------------------------------------
# Edit B: "Aplicación Web" paragraph - remove "o llevar", append
# " o materiales antes mencionados" after "medicamentos"
# ---------------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(5).Range
Write-Host "B para5 text before:" $p5.Text
$fragB = "<w:p $wns><w:pPr><w:pStyle w:val=`"Prrafodelista`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:jc w:val=`"both`"/></w:pPr>" + `
  "<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Aplicación Web (Angular):</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> Esta aplicación podrá ser accedida desde cualquier navegador de internet donde podrán visualizar un mapa </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">en tiempo real </w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`">con todos los lugares posibles para recoger </w:t></w:r>" + `
  "<w:r><w:t>medicamentos</w:t></w:r>" + `
  "<w:r><w:t xml:space=`"preserve`"> o materiales antes mencionados</w:t></w:r>" + `
  "<w:r><w:t>.</w:t></w:r>" + `
  "</w:p>"
$resB = $p5.InsertXML($fragB)
Write-Host "B result:" $resB
